$d = $word.ActiveDocument

# 1) Merge the split run "...paragraphe  {paragrapheAttestationConformit" + "e" + "}) ou un renouvellement..."
#    into a single contiguous run of text (no visible text change, just de-fragmenting the run boundaries
#    that Word naturally collapses when a Find/Replace rewrites across them).
$d.Content.Find.Execute(
    "paragrapheAttestationConformite}",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "paragrapheAttestationConformite}",
    2
) | Out-Null

# 2) Replace the "Copie : ..." line with the DREAL-templated version.
#    NB: the source text uses non-breaking spaces (U+00A0) before the French
#    punctuation (": ", "/ ", "; ") - preserve that exactly for the part of
#    the sentence that is untouched by the edit, only the "[DREAL concernée] ;
#    [CRE]" tail is actually being replaced.
$nbsp = [char]0x00A0
$findText    = "Copie${nbsp}: [EDF OA${nbsp}/ EDF SEI]${nbsp}; [DREAL concernée]${nbsp}; [CRE]"
$replaceText = "Copie${nbsp}: [EDF OA${nbsp}/ EDF SEI]${nbsp}; {#dreal}DREAL {dreal}{/dreal}{^dreal}DREAL concernée{/dreal} ; CRE"

$d.Content.Find.Execute(
    $findText,
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    $replaceText,
    2
) | Out-Null
